$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (10) shifting the old J:N (effect_type..is_parent)
# over to K:O, and set the new column's header.
$ws.Columns.Item(10).Insert()
$ws.Cells.Item(1, 10).Value = "resource_request_time_reduction"

# Best-effort column width for the newly inserted column (matches the
# bestFit width used by the other columns in this sheet).
$ws.Columns.Item(10).ColumnWidth = 37.705

# Add the new skill/quest row 21: "Speedy Resources"
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Speedy Resources"
$ws.Cells.Item(21, 3).Value = "When requesting resources from other kingdoms, the time to travel between kingdoms will reduce by 15% per level for a max of 60%."
$ws.Cells.Item(21, 4).Value = 4
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 10).Value = 0.15
$ws.Cells.Item(21, 11).Value = 14
$ws.Cells.Item(21, 12).Value = 17
$ws.Cells.Item(21, 13).Value = 1
$ws.Cells.Item(21, 14).Value = 1
$ws.Cells.Item(21, 15).Value = 1
